# TIME.xlsx reference workbook cleanup:
#   - Drop the unused, empty "Sheet1" tab.
#   - Rename the remaining "A  Regular Holidays" tab to "Timesheet"
#     (Excel automatically rewrites the ExternalData_1 defined name and
#     any sheet-qualified references to use the new tab name).

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Sheet1").Delete() | Out-Null
$wb.Worksheets.Item("A  Regular Holidays").Name = "Timesheet"
